$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text runs)
# ---------------------------------------------------------------------------
# A8 holds "Volume 29   Number  47" -> bump issue number to 48
$ws.Range("A8").Characters(21, 2).Text = "48"

# C9 holds "Report Covering the Week  11/21/2022  Through  11/27/2022"
# -> new reporting week
$ws.Range("C9").Characters(27, 10).Text = "11/28/2022"
$ws.Range("C9").Characters(48, 10).Text = "12/4/2022"

# ---------------------------------------------------------------------------
# Helper cells that keep a stable "text" representation we can clone styles
# and shared-string content from (row 22 stays untouched by this revision).
# ---------------------------------------------------------------------------
$zeroText = $ws.Range("C22")   # text "0"  (style 14)
$naText   = $ws.Range("E22")   # text "***.*" (style 14)

# ---------------------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------------------
$zeroText.Copy($ws.Range("C14"))
$ws.Range("F14").Value = 4
$ws.Range("H14").Value = 300
$ws.Range("M14").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 39
$ws.Range("K15").Value = -69.230769230769
$ws.Range("L15").Value = -65.714285714285
$ws.Range("N15").Value = -67.567567567567

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 167
$ws.Range("J16").Value = 132
$ws.Range("K16").Value = 26.515151515151
$ws.Range("L16").Value = 3.726708074534
$ws.Range("M16").Value = -47.484276729559
$ws.Range("N16").Value = -82.730093071354

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 20.833333333333
$ws.Range("I17").Value = 394
$ws.Range("J17").Value = 380
$ws.Range("K17").Value = 3.684210526315
$ws.Range("L17").Value = -2.475247524752
$ws.Range("M17").Value = 49.242424242424
$ws.Range("N17").Value = 1.546391752577

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -48
$ws.Range("I18").Value = 197
$ws.Range("J18").Value = 148
$ws.Range("K18").Value = 33.108108108108
$ws.Range("L18").Value = 1.546391752577
$ws.Range("M18").Value = -46.027397260274
$ws.Range("N18").Value = -87.107329842931

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -27.777777777777
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -15.873015873015
$ws.Range("I19").Value = 588
$ws.Range("J19").Value = 450
$ws.Range("K19").Value = 30.666666666666
$ws.Range("L19").Value = 6.137184115523
$ws.Range("M19").Value = 40.334128878281
$ws.Range("N19").Value = 6.137184115523

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 233.333333333333
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 21.052631578947
$ws.Range("I20").Value = 255
$ws.Range("J20").Value = 143
$ws.Range("K20").Value = 78.321678321678
$ws.Range("L20").Value = -7.272727272727
$ws.Range("M20").Value = -28.969359331476
$ws.Range("N20").Value = -92.102818209972

# ---------------------------------------------------------------------------
# Row 21 (bold totals row)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -10.810810810810
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 147
$ws.Range("H21").Value = -10.204081632653
$ws.Range("I21").Value = 1623
$ws.Range("J21").Value = 1297
$ws.Range("K21").Value = 25.134926754047
$ws.Range("L21").Value = -0.490496627835
$ws.Range("M21").Value = -7.679180887372
$ws.Range("N21").Value = -75.866171003717

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -41.935483870967
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = -37.142857142857
$ws.Range("I24").Value = 1395
$ws.Range("J24").Value = 994
$ws.Range("K24").Value = 40.342052313883
$ws.Range("L24").Value = 49.678111587982
$ws.Range("M24").Value = 83.552631578947

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -2.5
$ws.Range("I25").Value = 609
$ws.Range("J25").Value = 472
$ws.Range("K25").Value = 29.025423728813
$ws.Range("L25").Value = 24.032586558044
$ws.Range("M25").Value = 8.75

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = -44
$ws.Range("L26").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 60
$ws.Range("I27").Value = 60
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = 22.448979591836
$ws.Range("L27").Value = 22.448979591836

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$zeroText.Copy($ws.Range("C28"))
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -30.555555555555
$ws.Range("L28").Value = -28.571428571428
$ws.Range("M28").Value = -46.808510638297
$ws.Range("N28").Value = -67.948717948717

# ---------------------------------------------------------------------------
# Row 29
# ---------------------------------------------------------------------------
$zeroText.Copy($ws.Range("C29"))
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 2
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("J29").Value = 30
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -28.571428571428
$ws.Range("M29").Value = -35.483870967741
$ws.Range("N29").Value = -71.014492753623
